$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they are not
# auto-converted to numbers (preserves literal formatting, e.g. trailing zeros).
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14",
    "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25",
    "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35",
    "D36", "D37", "D39", "D40", "D41", "D42", "D44", "D45", "D47", "D48",
    "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values (Price + Volume(1h), and the reordered rows 48-50).
$ws.Range("D2").Value = "29.221.70"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.829.06"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "234.62"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").Value = "0.5974"
$ws.Range("E6").Value = "  -4.61%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "0.06962"
$ws.Range("E8").Value = "  -5.99%  "
$ws.Range("D9").Value = "0.2757"
$ws.Range("E9").Value = "  -4.64%  "
$ws.Range("D10").Value = "23.29"
$ws.Range("E10").Value = "  -6.09%  "
$ws.Range("D11").Value = "0.07621"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "1.826.45"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "4.769"
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").Value = "0.6217"
$ws.Range("E14").Value = "  -7.57%  "
$ws.Range("D15").Value = "0.000009608"
$ws.Range("E15").Value = "  -6.11%  "
$ws.Range("D16").Value = "78.46"
$ws.Range("E16").Value = "  -3.95%  "
$ws.Range("D17").Value = "29.179.39"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "5.751"
$ws.Range("E18").Value = "  -7.86%  "
$ws.Range("D19").Value = "222.37"
$ws.Range("E19").Value = "  -5.29%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "11.56"
$ws.Range("E21").Value = "  -6.00%  "
$ws.Range("D22").Value = "6.875"
$ws.Range("E22").Value = "  -5.76%  "
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").Value = "156.42"
$ws.Range("D25").Value = "0.1294"
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("D26").Value = "7.938"
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("D27").Value = "16.51"
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").Value = "0.06635"
$ws.Range("E28").Value = "  -7.58%  "
$ws.Range("D29").Value = "1.451"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("D30").Value = "1.439"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").Value = "3.825"
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("D32").Value = "3.751"
$ws.Range("E32").Value = "  -7.33%  "
$ws.Range("D33").Value = "1.094"
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("D34").Value = "1.717"
$ws.Range("E34").Value = "  -5.43%  "
$ws.Range("D35").Value = "0.6404"
$ws.Range("E35").Value = "  -7.82%  "
$ws.Range("D36").Value = "2.542"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "2.739"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Value = "1.199.05"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").Value = "0.01740"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("D40").Value = "6.455"
$ws.Range("E40").Value = "  -6.64%  "
$ws.Range("D41").Value = "0.8999"
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "1.982.44"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "100.03"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "62.08"
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("D47").Value = "8.477"
$ws.Range("E47").Value = "  -4.21%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.4555"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05507"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.565"
$ws.Range("E50").Value = "  -7.87%  "
$ws.Range("D51").Value = "6.364"
$ws.Range("E51").Value = "  -8.39%  "
